# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# New K values per row, calculated/regenerated from the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 4;
    3 = 5;
    4 = 6;
    5 = 6;
    6 = 6;
    7 = 9;
    8 = 8;
    9 = 8;
    10 = 8;
    11 = 2;
    12 = 10;
    13 = 5;
    14 = 11;
    15 = 9;
    16 = 8;
    17 = 5;
    18 = 8;
    19 = 7;
    20 = 10;
    21 = 7;
    22 = 2;
    23 = 8;
    24 = 10;
    25 = 8;
    26 = 6;
    27 = 6;
    28 = 5;
    29 = 10;
    30 = 9;
    31 = 8;
    32 = 6;
    33 = 7;
    34 = 9;
    35 = 9;
    36 = 4;
    37 = 5;
    38 = 7;
    39 = 9;
    40 = 3;
    41 = 3;
    42 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
